$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @("M1", "Il1a", "Il1r2", "FAPs", 3, 1, 4.072082666666667, 12.216248, 0.3979924983064649, 0.3979924983064649, 1, 0.3333333333333333, 0.05504633333333334, 0.165139, 0.002371469359192472, 0.002371469359192472, 0.2241532198302222, 2.017378978472, 0.0009438270149222433, 0.0009438270149222434)
$data[3] = @("M1", "Il1a", "Il1r2", "sCs", 3, 1, 4.072082666666667, 12.216248, 0.3979924983064649, 0.3979924983064649, 3, 1, 1.016579, 3.049737, 0.04379557735662424, 0.04379557735662425, 4.139593725197333, 37.256343526776, 0.01743031124693692, 0.01743031124693693)
$data[4] = @("M1", "Il1a", "Il1r2", "M1", 3, 1, 4.072082666666667, 12.216248, 0.3979924983064649, 0.3979924983064649, 3, 1, 8.671405, 26.014215, 0.3735756773139306, 0.3735756773139307, 35.31067799614667, 317.79610196532, 0.148680317120701, 0.148680317120701)
$data[5] = @("M1", "Il1a", "Il1r2", "M2", 3, 1, 4.072082666666667, 12.216248, 0.3979924983064649, 0.3979924983064649, 3, 1, 12.91794, 38.75382, 0.5565220613038737, 0.5565220613038738, 52.60291956304, 473.42627606736, 0.2214916055409923, 0.2214916055409923)
$data[6] = @("M1", "Il1a", "Il1r2", "ECs", 3, 1, 4.072082666666667, 12.216248, 0.3979924983064649, 0.3979924983064649, 3, 1, 0.5509396666666667, 1.652819, 0.02373521466637888, 0.02373521466637889, 2.243471867012444, 20.191246803112, 0.009446437382912379, 0.009446437382912379)
$data[7] = @("M2", "Il1a", "Il1r2", "FAPs", 3, 1, 6.159473666666667, 18.478421, 0.6020075016935351, 0.6020075016935351, 1, 0.3333333333333333, 0.05504633333333334, 0.165139, 0.002371469359192472, 0.002371469359192472, 0.3390564406132223, 3.051507965519, 0.001427642344270229, 0.001427642344270229)
$data[8] = @("M2", "Il1a", "Il1r2", "sCs", 3, 1, 6.159473666666667, 18.478421, 0.6020075016935351, 0.6020075016935351, 3, 1, 1.016579, 3.049737, 0.04379557735662424, 0.04379557735662425, 6.261591580586333, 56.354324225277, 0.02636526610968732, 0.02636526610968732)
$data[9] = @("M2", "Il1a", "Il1r2", "M1", 3, 1, 6.159473666666667, 18.478421, 0.6020075016935351, 0.6020075016935351, 3, 1, 8.671405, 26.014215, 0.3735756773139306, 0.3735756773139307, 53.41129075050167, 480.701616754515, 0.2248953601932296, 0.2248953601932296)
$data[10] = @("M2", "Il1a", "Il1r2", "M2", 3, 1, 6.159473666666667, 18.478421, 0.6020075016935351, 0.6020075016935351, 3, 1, 12.91794, 38.75382, 0.5565220613038737, 0.5565220613038738, 79.56771125758002, 716.1094013182201, 0.3350304557628814, 0.3350304557628815)
$data[11] = @("M2", "Il1a", "Il1r2", "ECs", 3, 1, 6.159473666666667, 18.478421, 0.6020075016935351, 0.6020075016935351, 3, 1, 0.5509396666666667, 1.652819, 0.02373521466637888, 0.02373521466637889, 3.393498368755445, 30.541485318799, 0.0142887772834665, 0.01428877728346651)

foreach ($r in 2..11) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $rowVals[$i]
    }
}
